# Atualizações após dados finais recebidos
# O professor enviou os últimos dados em 24/03/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 2
$ws.Range("F2").Value = 23
$ws.Range("G2").Value = 24
$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 36
$ws.Range("J2").Value = 37

# Row 3
$ws.Range("F3").Value = 22
$ws.Range("G3").Value = 23
$ws.Range("H3").Value = 34
$ws.Range("I3").Value = 35
$ws.Range("J3").Value = 36

# Row 4
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 24
$ws.Range("J4").Value = 25

# Row 5
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 23
$ws.Range("J5").Value = 24

# Row 6
$ws.Range("H6").Value = 11

# Row 7
$ws.Range("H7").Value = 10

# Update the active selection to match the saved view state
$ws.Range("P10").Select()
